$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "Descripción de Usuarios, Roles y Actores" ->
# "Descripción de Usuarios" + " y " + "R" + "oles " (4 separate runs,
# dropping ", Actores" from the title of that section heading).
#
# Plain Find&Replace (or InsertAfter chaining) collapses the inserted
# text back into a single run because the runtime coalesces adjacent runs
# that end up with identical (here: empty) run properties. Temporarily
# bracketing each insertion point with a bookmark prevents that
# coalescing; deleting the bookmarks afterwards leaves the four runs
# intact without leaving any bookmark behind.

$r1 = $d.Content
$found1 = $r1.Find.Execute("Descripción de Usuarios, Roles y Actores", $true, $false, $false, $false, $false, $true, 1, $false, "Descripción de Usuarios", 2)

if ($found1) {
    $r1.Collapse(0)
    $d.Bookmarks.Add("ZZ_tmp_split_1", $r1)
    $r1.InsertAfter(" y ")

    $r1.Collapse(0)
    $d.Bookmarks.Add("ZZ_tmp_split_2", $r1)
    $r1.InsertAfter("R")

    $r1.Collapse(0)
    $d.Bookmarks.Add("ZZ_tmp_split_3", $r1)
    $r1.InsertAfter("oles ")

    $d.Bookmarks("ZZ_tmp_split_1").Delete()
    $d.Bookmarks("ZZ_tmp_split_2").Delete()
    $d.Bookmarks("ZZ_tmp_split_3").Delete()
}

# --- Change 2 -----------------------------------------------------------
# " (" + "Backend" + " y " (the middle run wrapped in spell-check
# proofErr markers) is merged back into a single run " (Backend y ",
# removing the now-unneeded proofErr markers around "Backend".

$r2 = $d.Content
$null = $r2.Find.Execute(" (Backend y ", $true, $false, $false, $false, $false, $true, 1, $false, " (Backend y ", 2)
